$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.5589887640449438
$ws1.Range("C2").Value = 0.5314057826520439
$ws1.Range("E2").Value = 0.69355888093689
$ws1.Range("F2").Value = 0.8489964956992673
$ws1.Range("G2").Value = 0.9655124364244408
$ws1.Range("H2").Value = 0.7962764241327553
$ws1.Range("J2").Value = 470
$ws1.Range("K2").Value = 64

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")
$ws2.Range("B2").Value = 0.9846153846153847
$ws2.Range("C2").Value = 0.1198501872659176
$ws2.Range("D2").Value = 0.2136894824707846

$ws2.Range("B3").Value = 0.5314057826520439
$ws2.Range("D3").Value = 0.69355888093689

$ws2.Range("B4").Value = 0.5589887640449438
$ws2.Range("C4").Value = 0.5589887640449438
$ws2.Range("D4").Value = 0.5589887640449438
$ws2.Range("E4").Value = 0.5589887640449438

$ws2.Range("B5").Value = 0.7580105836337143
$ws2.Range("C5").Value = 0.5589887640449438
$ws2.Range("D5").Value = 0.4536241817038373

$ws2.Range("B6").Value = 0.7580105836337143
$ws2.Range("C6").Value = 0.5589887640449438
$ws2.Range("D6").Value = 0.4536241817038373

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")
$ws3.Range("B2").Value = 64
$ws3.Range("C2").Value = 470
